# Update column F ("dSF") values for the rows that changed on repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = 4
$ws.Range("F6").Value  = 4
$ws.Range("F10").Value = 0
$ws.Range("F18").Value = -6
$ws.Range("F19").Value = -8
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = -9
$ws.Range("F24").Value = -1
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("F29").Value = -10
